{"js": "// Strengthen \"What Not To Do\" warnings (v0.5)\n// Replace the 12 ListBullet paragraph texts under \"WHAT NOT TO DO (YEARS OF MISTAKES)\"\n// with the expanded/strengthened wording from the diff.\n\nconst replacements = [\n  [\n    \"Avoid **starting with 15+ crop types**.\",\n    \"Don\\u2019t run **15+ crop types** in year one. It fragments attention, guarantees missed tasks, and produces a season of mediocre output.\"\n  ],\n  [\n    \"Don\\u2019t use **tiny containers** for anything you expect to eat weekly.\",\n    \"Don\\u2019t bet on **small containers** for weekly food. They dry fast, swing temperature hard, and quietly cap yield no matter how \\u201cgood\\u201d the plant looks.\"\n  ],\n  [\n    \"Avoid **watering by mood**.\",\n    \"Don\\u2019t **water opportunistically**. Irregular watering triggers blossom drop, cracked fruit, bitter greens, and weeks of lost momentum.\"\n  ],\n  [\n    \"Don\\u2019t treat fertilizer as a substitute for stable water.\",\n    \"Don\\u2019t use **fertilizer to compensate for water stress**. It turns stress into pest pressure and weak growth.\"\n  ],\n  [\n    \"Avoid **late pest response**.\",\n    \"Don\\u2019t wait for visible damage before acting on pests. By the time leaves look chewed, you\\u2019ve already paid in yield.\"\n  ],\n  [\n    \"Don\\u2019t plant calorie crops in shade and blame the crop.\",\n    \"Don\\u2019t put calorie crops into shade and pretend effort will replace photons. The harvest will never match the space.\"\n  ],\n  [\n    \"Avoid **one-time planting** of greens.\",\n    \"Don\\u2019t do a single planting of greens and call it a system. You get a short peak, then empty containers.\"\n  ],\n  [\n    \"Don\\u2019t keep weak plants \\u201cto see if they recover.\\u201d\",\n    \"Don\\u2019t keep weak plants \\u201cto see if they recover.\\u201d In small spaces, a stalled plant is a space tax you pay every day.\"\n  ],\n  [\n    \"Avoid variety obsession before water and soil are stable.\",\n    \"Don\\u2019t obsess over varieties while water and soil are unstable. You\\u2019ll blame genetics for operational failure.\"\n  ],\n  [\n    \"Don\\u2019t ignore wind; mitigate or relocate containers.\",\n    \"Don\\u2019t ignore wind exposure. Wind-driven drying and stem stress will cut output and break plants at the worst time.\"\n  ],\n  [\n    \"Avoid unknown compost/inputs in containers.\",\n    \"Don\\u2019t add unknown compost/inputs to containers. One contaminated or too-hot input can stunt everything for the season.\"\n  ],\n  [\n    \"Don\\u2019t scale plant count until you can run the system through the hottest month.\",\n    \"Don\\u2019t scale plant count until you\\u2019ve run your setup through the hottest month without repeated wilting or missed waterings.\"\n  ]\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Build a lookup of old text -> new text for fast matching.\nconst map = new Map(replacements);\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  const current = para.text;\n  if (map.has(current)) {\n    para.insertText(map.get(current), \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Strengthen \"What Not To Do\" warnings (v0.5)\n# Replace the 12 ListBullet paragraph texts under \"WHAT NOT TO DO (YEARS OF MISTAKES)\"\n# with the expanded/strengthened wording, matched via Find.Execute on the exact old text.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n    ,@('Avoid **starting with 15+ crop types**.', 'Don\u2019t run **15+ crop types** in year one. It fragments attention, guarantees missed tasks, and produces a season of mediocre output.')\n    ,@('Don\u2019t use **tiny containers** for anything you expect to eat weekly.', 'Don\u2019t bet on **small containers** for weekly food. They dry fast, swing temperature hard, and quietly cap yield no matter how \u201cgood\u201d the plant looks.')\n    ,@('Avoid **watering by mood**.', 'Don\u2019t **water opportunistically**. Irregular watering triggers blossom drop, cracked fruit, bitter greens, and weeks of lost momentum.')\n    ,@('Don\u2019t treat fertilizer as a substitute for stable water.', 'Don\u2019t use **fertilizer to compensate for water stress**. It turns stress into pest pressure and weak growth.')\n    ,@('Avoid **late pest response**.', 'Don\u2019t wait for visible damage before acting on pests. By the time leaves look chewed, you\u2019ve already paid in yield.')\n    ,@('Don\u2019t plant calorie crops in shade and blame the crop.', 'Don\u2019t put calorie crops into shade and pretend effort will replace photons. The harvest will never match the space.')\n    ,@('Avoid **one-time planting** of greens.', 'Don\u2019t do a single planting of greens and call it a system. You get a short peak, then empty containers.')\n    ,@('Don\u2019t keep weak plants \u201cto see if they recover.\u201d', 'Don\u2019t keep weak plants \u201cto see if they recover.\u201d In small spaces, a stalled plant is a space tax you pay every day.')\n    ,@('Avoid variety obsession before water and soil are stable.', 'Don\u2019t obsess over varieties while water and soil are unstable. You\u2019ll blame genetics for operational failure.')\n    ,@('Don\u2019t ignore wind; mitigate or relocate containers.', 'Don\u2019t ignore wind exposure. Wind-driven drying and stem stress will cut output and break plants at the worst time.')\n    ,@('Avoid unknown compost/inputs in containers.', 'Don\u2019t add unknown compost/inputs to containers. One contaminated or too-hot input can stunt everything for the season.')\n    ,@('Don\u2019t scale plant count until you can run the system through the hottest month.', 'Don\u2019t scale plant count until you\u2019ve run your setup through the hottest month without repeated wilting or missed waterings.')\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $found = $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        Write-Output (\"NOT FOUND: \" + $oldText)\n    }\n}\n\n"}
